$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Clear the whole used range and start the content layout fresh so that no
#    stale formatting (alignment, old number formats, etc.) leaks through.
# ---------------------------------------------------------------------------
$ws.Range("A1:B26").ClearFormats()
$ws.Range("A1:B26").ClearContents()

# ---------------------------------------------------------------------------
# 2. Header row (row 1)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Burnout Cost Factors"
$ws.Range("B1").Value = "Inputs"

# ---------------------------------------------------------------------------
# 3. Input labels (column A) & blank input cells (column B), rows 2-13
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Number of Employees"
$ws.Range("A3").Value = "Turnover Rate Due to Burnout (%)"
$ws.Range("A4").Value = "Average Salary per Employee ($)"
$ws.Range("A5").Value = "Replacement Cost Factor (as % of Salary)"
$ws.Range("A6").Value = "Productivity Loss Due to Burnout (%)"
$ws.Range("A7").Value = "Additional Sick Days per Burnt-out Employee"
$ws.Range("A8").Value = "Average Mental Health Claim per Employee ($)"
$ws.Range("A9").Value = "Manager Time Lost to Burnout Issues (hrs/month)"
$ws.Range("A10").Value = "Average Manager Hourly Rate ($)"
$ws.Range("A11").Value = "Opportunity Cost (Missed Revenue or Projects $)"
$ws.Range("A12").Value = "Presenteeism Cost Factor (as % of Salary)"
$ws.Range("A13").Value = "Burnout Reduction Target (%)"
# Row 14 left blank as a spacer row (both columns)

# ---------------------------------------------------------------------------
# 4. Second section header (row 15)
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "Calculated Costs"
$ws.Range("B15").Value = "Results"

# ---------------------------------------------------------------------------
# 5. Result labels (column A) & formulas (column B), rows 16-26
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = "Turnover Cost ($)"
$ws.Range("B16").Formula = "=(B3*B4)*(B2/100)"

$ws.Range("A17").Value = "Productivity Loss Cost Due to Burnout ($)"
$ws.Range("B17").Formula = "=B3*B2*(B5/100)"

$ws.Range("A18").Value = "Sick Leave Cost ($)"
$ws.Range("B18").Formula = "=B2*B6*(B4/260)"

$ws.Range("A19").Value = "Mental Health Claim Cost ($)"
$ws.Range("B19").Formula = "=B2*B7"

$ws.Range("A20").Value = "Manager Time Cost ($)"
$ws.Range("B20").Formula = "=B8*B9*12"

$ws.Range("A21").Value = "Presenteeism Cost ($)"
$ws.Range("B21").Formula = "=B3*B2*(B11/100)"

$ws.Range("A22").Value = "Opportunity Cost ($)"
$ws.Range("B22").Formula = "=B10"

$ws.Range("A23").Value = "Total Burnout Cost ($)"
$ws.Range("B23").Formula = "=SUM(B16:B22)"

$ws.Range("A24").Value = "Projected Savings if Burnout Reduced ($)"
$ws.Range("B24").Formula = "=B23*(B12/100)"

$ws.Range("A25").Value = "Cost per Employee ($)"
$ws.Range("B25").Formula = "=B23/B2"

$ws.Range("A26").Value = "% Payroll Lost to Burnout"
$ws.Range("B26").Formula = "=B23/(B3*B4)"

# ---------------------------------------------------------------------------
# 6. Formatting
# ---------------------------------------------------------------------------

# -- Section header rows (1 and 15): bold white text on blue fill, thin border
foreach ($addr in @("A1:B1", "A15:B15")) {
    $headerRange = $ws.Range($addr)
    $headerRange.Font.Bold = $true
    $headerRange.Font.Color = 16777215
    $headerRange.Interior.Color = 12874308
    $headerRange.Borders.LineStyle = 1
}

# -- Label column (A2:A14 and A16:A26): bold text on light-blue fill, thin border
$labelRange = $ws.Range("A2:A14")
$labelRange.Font.Bold = $true
$labelRange.Interior.Color = 15917529
$labelRange.Borders.LineStyle = 1

$labelRange2 = $ws.Range("A16:A26")
$labelRange2.Font.Bold = $true
$labelRange2.Interior.Color = 15917529
$labelRange2.Borders.LineStyle = 1

# -- Input cells (B2:B14): plain number format, thin border
$inputRange = $ws.Range("B2:B14")
$inputRange.NumberFormat = "#,##0"
$inputRange.Borders.LineStyle = 1

# -- Result cells (B16:B25): currency format, thin border
$resultRange = $ws.Range("B16:B25")
$resultRange.NumberFormat = "\$#,##0"
$resultRange.Borders.LineStyle = 1

# -- Final result cell (B26): percentage format, thin border
$ws.Range("B26").NumberFormat = "0.00%"
$ws.Range("B26").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 7. Column widths
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 44.83
$ws.Columns("B").ColumnWidth = 24.83

# ---------------------------------------------------------------------------
# 8. Selection
# ---------------------------------------------------------------------------
$ws.Range("B2:B4").Select()
